# SAV-MVH.xlsx update: "Actualizando HU y MVH"
# Renumbers / re-labels several HU (historia de usuario) entries in the
# RELEASE 3 / RELEASE 4 block (rows 11-18) and refreshes the active
# selection in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 11: "H23 - Visualizar tareas Asignadas" -> "H26 - Visualizar tareas
# Asignadas" (and it moves from A11/B11 merge-less pair into just A11).
# B11 ("H28 - Generar examen") and G11 ("H29 - Visualizar reporte de
# alumnos") move down into the new row 17, so clear them here.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "H26 - Visualizar tareas Asignadas"
$ws.Range("B11").Clear()
$ws.Range("G11").ClearContents()
$ws.Range("I11").Value = "H28 - Matricular estudiantes"

# ---------------------------------------------------------------------
# Row 12: "H26 - Agregar tareas" -> "H27 - Agregar tareas"
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "H27 - Agregar tareas"

# ---------------------------------------------------------------------
# Row 13: B13 picks up the bold "H##" header look (style used by B12/G12),
# H13 "H 24 - Editar estudiante" -> "H24 - Editar estudiante"
# ---------------------------------------------------------------------
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("H13").Value = "H24 - Editar estudiante"

# ---------------------------------------------------------------------
# Row 14: B14/C14 adopt the same look as B13/C13 respectively, and H14
# gets the new "H25 - Eliminar estudiante" entry.
# ---------------------------------------------------------------------
$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("H14").Value = "H25 - Eliminar estudiante"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 16 becomes a thin gray separator row, matching rows 2/7/10.
# ---------------------------------------------------------------------
$ws.Range("A2:J2").Copy()
$ws.Range("A16:J16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows(16).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# Row 17 becomes a brand-new content row: "H29 - Generar examen" (B17),
# "H30 - Visualizar reporte de alumnos" (G17) and the new "RELEASE 4"
# label (J17). A17/F17 (previously blank, just fill) are cleared out.
# ---------------------------------------------------------------------
$ws.Range("A17").Clear()
$ws.Range("F17").Clear()

$ws.Range("B12").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "H29 - Generar examen"

$ws.Range("B12").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value = "H30 - Visualizar reporte de alumnos"

$ws.Range("J15").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("J17").Value = "RELEASE 4"
$excel.CutCopyMode = 0

$ws.Rows(17).RowHeight = 28.8

# ---------------------------------------------------------------------
# Row 18 loses its leading/trailing blank cells (A18, F18) and grows a
# touch taller.
# ---------------------------------------------------------------------
$ws.Range("A18").Clear()
$ws.Range("F18").Clear()
$ws.Rows(18).RowHeight = 31.2

# ---------------------------------------------------------------------
# Row 24 reverts to the default (non-custom) row height, row 26 becomes
# a fixed-height row, matching the re-shuffled layout below.
# ---------------------------------------------------------------------
$ws.Rows(24).EntireRow.AutoFit()
$ws.Rows(26).RowHeight = 14.4

# ---------------------------------------------------------------------
# Refresh the view: scroll so row 4 is at the top and select C13 (as
# left behind by whoever made this edit).
# ---------------------------------------------------------------------
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 4
